$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C30").Value = 8789
$ws.Range("D30").Value = 8789
$ws.Range("E30").Value = 8789
$ws.Range("I30").Value = 71.406227644719095
$ws.Range("J30").Value = 58.367133751517699
$ws.Range("K30").Value = 152.51930729799

$ws.Range("C31").Value = 8789
$ws.Range("D31").Value = 8789
$ws.Range("E31").Value = 8789
$ws.Range("I31").Value = 405.34160708664899
$ws.Range("J31").Value = 465.57881909953397
$ws.Range("K31").Value = 543.86826505647298

$ws.Range("C32").Value = 8789
$ws.Range("D32").Value = 8789
$ws.Range("E32").Value = 8789
$ws.Range("I32").Value = 677.077774916301
$ws.Range("J32").Value = 811.86554716529895
$ws.Range("K32").Value = 948.41402235630198

$ws.Range("K33").Select()
